$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns H (AgTests) and I (AgPosit) for rows 303-335 with corrected values
$ws.Range("H303").Value = 9876
$ws.Range("I303").Value = 632
$ws.Range("H304").Value = 5690
$ws.Range("I304").Value = 479
$ws.Range("H305").Value = 3511
$ws.Range("I305").Value = 301
$ws.Range("H306").Value = 70451
$ws.Range("I306").Value = 7164
$ws.Range("H307").Value = 75237
$ws.Range("I307").Value = 6440
$ws.Range("H308").Value = 15750
$ws.Range("I308").Value = 1095
$ws.Range("H309").Value = 74302
$ws.Range("I309").Value = 5286
$ws.Range("H310").Value = 74236
$ws.Range("H311").Value = 63019
$ws.Range("H313").Value = 61401
$ws.Range("I313").Value = 3144
$ws.Range("H314").Value = 63431
$ws.Range("I314").Value = 3143
$ws.Range("H315").Value = 65804
$ws.Range("H316").Value = 49160
$ws.Range("H317").Value = 61385
$ws.Range("H319").Value = 41184
$ws.Range("I319").Value = 1628
$ws.Range("H320").Value = 76616
$ws.Range("H321").Value = 90377
$ws.Range("H322").Value = 106212
$ws.Range("I323").Value = 2304
$ws.Range("H324").Value = 230076
$ws.Range("I324").Value = 2635
$ws.Range("H325").Value = 700175
$ws.Range("I325").Value = 5771
$ws.Range("H326").Value = 417169
$ws.Range("I326").Value = 3686
$ws.Range("H327").Value = 235476
$ws.Range("I327").Value = 2873
$ws.Range("I329").Value = 1718
$ws.Range("H331").Value = 147127
$ws.Range("I331").Value = 2543
$ws.Range("H332").Value = 408753
$ws.Range("I332").Value = 3991
$ws.Range("H333").Value = 249927
$ws.Range("I333").Value = 2658
$ws.Range("H334").Value = 197811
$ws.Range("I334").Value = 3283
$ws.Range("H335").Value = 121310
$ws.Range("I335").Value = 2796

# Append new row 336 (04.02.2021 data) - new file structure from 07.Feb.2021
$ws.Range("A336").Value = 44230
$ws.Range("B336").Value = 256903
$ws.Range("C336").Value = 230456
$ws.Range("D336").Value = 21471
$ws.Range("E336").Value = 11223
$ws.Range("F336").Value = 2077
$ws.Range("G336").Value = 4976
$ws.Range("H336").Value = 89413
$ws.Range("I336").Value = 2879
